$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1542.2
$ws.Range("I19").Value = 2221.4167
$ws.Range("J19").Value = 523.375
$ws.Range("K19").Value = 2221.4167
$ws.Range("L19").Value = 523.375
$ws.Range("M19").Value = -2046.4167
$ws.Range("N19").Value = -873.375
$ws.Range("H132").Value = 2402.1738
$ws.Range("I132").Value = 2666.75
$ws.Range("J132").Value = 1797.4286
$ws.Range("K132").Value = 8000.25
$ws.Range("L132").Value = 5392.2858
$ws.Range("M132").Value = -5470.25
$ws.Range("N132").Value = -10452.2858
$ws.Range("H137").Value = 21015.762
$ws.Range("I137").Value = 11195.667
$ws.Range("J137").Value = 24943.8
$ws.Range("K137").Value = 33587.001
$ws.Range("L137").Value = 74831.39999999999
$ws.Range("M137").Value = -31037.001
$ws.Range("N137").Value = -79931.39999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7757.2
$ws.Range("I61").Value = 4438.579
$ws.Range("J61").Value = 18266.166
$ws.Range("K61").Value = 4438.579
$ws.Range("L61").Value = 18266.166
$ws.Range("M61").Value = -4226.579
$ws.Range("N61").Value = -18690.166
$ws.Range("H74").Value = 9619.075000000001
$ws.Range("I74").Value = 10639.625
$ws.Range("J74").Value = 5536.875
$ws.Range("K74").Value = 10639.625
$ws.Range("L74").Value = 5536.875
$ws.Range("M74").Value = -9765.625
$ws.Range("N74").Value = -7284.875
$ws.Range("H77").Value = 9619.075000000001
$ws.Range("I77").Value = 10639.625
$ws.Range("J77").Value = 5536.875
$ws.Range("K77").Value = 53198.125
$ws.Range("L77").Value = 27684.375
$ws.Range("M77").Value = -48830.125
$ws.Range("N77").Value = -36420.375
$ws.Range("H132").Value = 4823.1313
$ws.Range("I132").Value = 3889.9707
$ws.Range("J132").Value = 12755
$ws.Range("K132").Value = 11669.9121
$ws.Range("L132").Value = 38265
$ws.Range("M132").Value = -9139.9121
$ws.Range("N132").Value = -43325
$ws.Range("H136").Value = 7757.2
$ws.Range("I136").Value = 4438.579
$ws.Range("J136").Value = 18266.166
$ws.Range("K136").Value = 13315.737
$ws.Range("L136").Value = 54798.49800000001
$ws.Range("M136").Value = -10765.737
$ws.Range("N136").Value = -59898.49800000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2375.9167
$ws.Range("I99").Value = 2396.4546
$ws.Range("J99").Value = 2150
$ws.Range("K99").Value = 2396.4546
$ws.Range("L99").Value = 2150
$ws.Range("M99").Value = -898.4546
$ws.Range("N99").Value = -5146
$ws.Range("H134").Value = 12823.808
$ws.Range("I134").Value = 6801.7646
$ws.Range("J134").Value = 24198.777
$ws.Range("K134").Value = 20405.2938
$ws.Range("L134").Value = 72596.33099999999
$ws.Range("M134").Value = -17870.2938
$ws.Range("N134").Value = -77666.33099999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4809.5454
$ws.Range("I31").Value = 3871.6
$ws.Range("J31").Value = 5591.1665
$ws.Range("K31").Value = 3871.6
$ws.Range("L31").Value = 5591.1665
$ws.Range("M31").Value = -3576.6
$ws.Range("N31").Value = -6181.1665
$ws.Range("H34").Value = 4809.5454
$ws.Range("I34").Value = 3871.6
$ws.Range("J34").Value = 5591.1665
$ws.Range("K34").Value = 3871.6
$ws.Range("L34").Value = 5591.1665
$ws.Range("M34").Value = -3669.6
$ws.Range("N34").Value = -5995.1665
$ws.Range("H58").Value = 4874.4287
$ws.Range("I58").Value = 2400.7
$ws.Range("J58").Value = 11058.75
$ws.Range("K58").Value = 2400.7
$ws.Range("L58").Value = 11058.75
$ws.Range("M58").Value = -2197.7
$ws.Range("N58").Value = -11464.75
$ws.Range("H99").Value = 8307.082
$ws.Range("I99").Value = 3797
$ws.Range("K99").Value = 3797
$ws.Range("M99").Value = -2299
$ws.Range("H126").Value = 8307.082
$ws.Range("I126").Value = 3797
$ws.Range("K126").Value = 11391
$ws.Range("M126").Value = -8921
$ws.Range("H132").Value = 19874.746
$ws.Range("I132").Value = 11855.25
$ws.Range("J132").Value = 39400.477
$ws.Range("K132").Value = 35565.75
$ws.Range("L132").Value = 118201.431
$ws.Range("M132").Value = -33035.75
$ws.Range("N132").Value = -123261.431
$ws.Range("H134").Value = 3766.5
$ws.Range("I134").Value = 3642.6
$ws.Range("J134").Value = 4014.3
$ws.Range("K134").Value = 10927.8
$ws.Range("L134").Value = 12042.9
$ws.Range("M134").Value = -8392.799999999999
$ws.Range("N134").Value = -17112.9
$ws.Range("H136").Value = 4874.4287
$ws.Range("I136").Value = 2400.7
$ws.Range("J136").Value = 11058.75
$ws.Range("K136").Value = 7202.099999999999
$ws.Range("L136").Value = 33176.25
$ws.Range("M136").Value = -4652.099999999999
$ws.Range("N136").Value = -38276.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1806.6111
$ws.Range("I5").Value = 749
$ws.Range("J5").Value = 1868.8235
$ws.Range("K5").Value = 2247
$ws.Range("L5").Value = 5606.470499999999
$ws.Range("M5").Value = -2135
$ws.Range("N5").Value = -5830.470499999999
$ws.Range("H107").Value = 725.6
$ws.Range("I107").Value = 611.1818
$ws.Range("J107").Value = 865.44446
$ws.Range("K107").Value = 1833.5454
$ws.Range("L107").Value = 2596.33338
$ws.Range("M107").Value = 86.45460000000003
$ws.Range("N107").Value = -6436.33338
$ws.Range("H135").Value = 1806.6111
$ws.Range("I135").Value = 749
$ws.Range("J135").Value = 1868.8235
$ws.Range("K135").Value = 6741
$ws.Range("L135").Value = 16819.4115
$ws.Range("M135").Value = -4206
$ws.Range("N135").Value = -21889.4115
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2900.8572
$ws.Range("I80").Value = 1462.2
$ws.Range("J80").Value = 6497.5
$ws.Range("K80").Value = 1462.2
$ws.Range("L80").Value = 6497.5
$ws.Range("M80").Value = -464.2
$ws.Range("N80").Value = -8493.5
$ws.Range("H83").Value = 2900.8572
$ws.Range("I83").Value = 1462.2
$ws.Range("J83").Value = 6497.5
$ws.Range("K83").Value = 7311
$ws.Range("L83").Value = 32487.5
$ws.Range("M83").Value = -2319
$ws.Range("N83").Value = -42471.5
$ws.Range("H102").Value = 2796
$ws.Range("I102").Value = 2853.6667
$ws.Range("K102").Value = 2853.6667
$ws.Range("M102").Value = -1231.6667
$ws.Range("H132").Value = 16485.77
$ws.Range("I132").Value = 14789.5625
$ws.Range("J132").Value = 19199.7
$ws.Range("K132").Value = 44368.6875
$ws.Range("L132").Value = 57599.10000000001
$ws.Range("M132").Value = -41838.6875
$ws.Range("N132").Value = -62659.10000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3988.1428
$ws.Range("I122").Value = 3050.8572
$ws.Range("J122").Value = 6331.357
$ws.Range("K122").Value = 9152.571599999999
$ws.Range("L122").Value = 18994.071
$ws.Range("M122").Value = -6702.571599999999
$ws.Range("N122").Value = -23894.071
$ws.Range("H132").Value = 12595.258
$ws.Range("I132").Value = 7565.761
$ws.Range("J132").Value = 25555.885
$ws.Range("K132").Value = 22697.283
$ws.Range("L132").Value = 76667.655
$ws.Range("M132").Value = -20167.283
$ws.Range("N132").Value = -81727.655
